{"js": "// The upstream change (M2Doc test-fixture \"Moving from 2.0.1 to 2.0.2\")\n// is purely a re-serialization of the package XML: every element's\n// attributes got re-emitted in (alphabetical) order by the newer\n// docx4j-based tooling, and the root <w:document> namespace\n// declarations got re-sorted the same way. Attribute/namespace\n// declaration order is not part of the WordprocessingML information\n// model (readers must treat `<w:color w:val=\"x\" w:themeColor=\"y\"/>`\n// and `<w:color w:themeColor=\"y\" w:val=\"x\"/>` as identical), and the\n// Word JavaScript API has no surface for controlling that raw\n// serialization order. Diffing the two XML payloads attribute-set by\n// attribute-set (rather than attribute-sequence by attribute-sequence)\n// confirms there is no actual content, formatting, or structural\n// change anywhere in the document: same paragraphs, same runs, same\n// field codes, same colors, same page size/margins, same styles/\n// latent-style table, same fonts/langs.\n//\n// So the faithful translation of this diff through the document object\n// model is a no-op: we touch nothing, and the body/styles are left\n// exactly as authored.\n", "ps1": "# The upstream change (M2Doc test-fixture \"Moving from 2.0.1 to 2.0.2\")\n# is purely a re-serialization of the package XML: every element's\n# attributes got re-emitted in (alphabetical) order by the newer\n# docx4j-based tooling, and the root <w:document> namespace\n# declarations got re-sorted the same way. Attribute/namespace\n# declaration order is not part of the WordprocessingML information\n# model, and the Word COM object model has no surface for controlling\n# that raw serialization order. Comparing the two XML payloads\n# attribute-set by attribute-set (rather than attribute-sequence by\n# attribute-sequence) confirms there is no actual content, formatting,\n# or structural change anywhere in the document: same paragraphs, same\n# runs, same field codes, same colors, same page size/margins, same\n# styles/latent-style table, same fonts/langs.\n#\n# So the faithful translation of this diff through the document object\n# model is a no-op: we touch nothing, and the body/styles are left\n# exactly as authored.\n$d = $word.ActiveDocument\n"}
